{"js": "// Update the worksheet date and the 25 three-digit \u00f7 one-digit division\n// problems/answers to the next day's generated set.\nconst replacements = [\n  [\"2025-04-05 Saturday\", \"2025-04-06 Sunday\"],\n  [\"221\u00f77=31, 4\", \"255\u00f77=36, 3\"],\n  [\"552\u00f75=110, 2\", \"287\u00f78=35, 7\"],\n  [\"285\u00f75=57, 0\", \"487\u00f77=69, 4\"],\n  [\"993\u00f79=110, 3\", \"507\u00f75=101, 2\"],\n  [\"798\u00f76=133, 0\", \"445\u00f73=148, 1\"],\n  [\"843\u00f72=421, 1\", \"475\u00f75=95, 0\"],\n  [\"355\u00f79=39, 4\", \"167\u00f76=27, 5\"],\n  [\"561\u00f73=187, 0\", \"427\u00f72=213, 1\"],\n  [\"647\u00f73=215, 2\", \"978\u00f72=489, 0\"],\n  [\"890\u00f75=178, 0\", \"517\u00f77=73, 6\"],\n  [\"292\u00f73=97, 1\", \"283\u00f78=35, 3\"],\n  [\"770\u00f73=256, 2\", \"158\u00f75=31, 3\"],\n  [\"283\u00f72=141, 1\", \"322\u00f73=107, 1\"],\n  [\"231\u00f79=25, 6\", \"763\u00f72=381, 1\"],\n  [\"837\u00f75=167, 2\", \"768\u00f77=109, 5\"],\n  [\"290\u00f72=145, 0\", \"976\u00f74=244, 0\"],\n  [\"807\u00f72=403, 1\", \"914\u00f76=152, 2\"],\n  [\"773\u00f73=257, 2\", \"930\u00f79=103, 3\"],\n  [\"690\u00f77=98, 4\", \"198\u00f73=66, 0\"],\n  [\"274\u00f74=68, 2\", \"702\u00f75=140, 2\"],\n  [\"476\u00f73=158, 2\", \"430\u00f78=53, 6\"],\n  [\"867\u00f72=433, 1\", \"702\u00f78=87, 6\"],\n  [\"268\u00f79=29, 7\", \"418\u00f77=59, 5\"],\n  [\"885\u00f74=221, 1\", \"869\u00f73=289, 2\"],\n  [\"684\u00f77=97, 5\", \"663\u00f75=132, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 three-digit \u00f7 one-digit division\n# problems/answers to the next day's generated set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-05 Saturday\", \"2025-04-06 Sunday\"),\n    @(\"221\u00f77=31, 4\", \"255\u00f77=36, 3\"),\n    @(\"552\u00f75=110, 2\", \"287\u00f78=35, 7\"),\n    @(\"285\u00f75=57, 0\", \"487\u00f77=69, 4\"),\n    @(\"993\u00f79=110, 3\", \"507\u00f75=101, 2\"),\n    @(\"798\u00f76=133, 0\", \"445\u00f73=148, 1\"),\n    @(\"843\u00f72=421, 1\", \"475\u00f75=95, 0\"),\n    @(\"355\u00f79=39, 4\", \"167\u00f76=27, 5\"),\n    @(\"561\u00f73=187, 0\", \"427\u00f72=213, 1\"),\n    @(\"647\u00f73=215, 2\", \"978\u00f72=489, 0\"),\n    @(\"890\u00f75=178, 0\", \"517\u00f77=73, 6\"),\n    @(\"292\u00f73=97, 1\", \"283\u00f78=35, 3\"),\n    @(\"770\u00f73=256, 2\", \"158\u00f75=31, 3\"),\n    @(\"283\u00f72=141, 1\", \"322\u00f73=107, 1\"),\n    @(\"231\u00f79=25, 6\", \"763\u00f72=381, 1\"),\n    @(\"837\u00f75=167, 2\", \"768\u00f77=109, 5\"),\n    @(\"290\u00f72=145, 0\", \"976\u00f74=244, 0\"),\n    @(\"807\u00f72=403, 1\", \"914\u00f76=152, 2\"),\n    @(\"773\u00f73=257, 2\", \"930\u00f79=103, 3\"),\n    @(\"690\u00f77=98, 4\", \"198\u00f73=66, 0\"),\n    @(\"274\u00f74=68, 2\", \"702\u00f75=140, 2\"),\n    @(\"476\u00f73=158, 2\", \"430\u00f78=53, 6\"),\n    @(\"867\u00f72=433, 1\", \"702\u00f78=87, 6\"),\n    @(\"268\u00f79=29, 7\", \"418\u00f77=59, 5\"),\n    @(\"885\u00f74=221, 1\", \"869\u00f73=289, 2\"),\n    @(\"684\u00f77=97, 5\", \"663\u00f75=132, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
